$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44308, 1, 4, 95.30617107457708),
    @(44309, 2, 5, 119.1327138432213),
    @(44310, 1, 5, 119.1327138432213),
    @(44311, 6, 10, 238.2654276864427),
    @(44312, 0, 10, 238.2654276864427)
)

$lastRow = 233
$startRow = $lastRow + 1

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Copy the formatting of the last existing row's cell A so the new
    # date cell picks up the same style (border/alignment/number format)
    # without registering a brand new style entry.
    $ws.Cells.Item($lastRow, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $row[0]

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
